$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.970.85"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.112.32"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'526.29"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").Value = "'142.29"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.107.84"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'0.441"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "'7.22"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("D13").Value = "3.636.58"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "'25.67"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "57.966.52"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "3.094.60"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'12.86"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'8.03"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "'341.87"
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("D25").Value = "'67.37"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'6.51"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'7.26"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "'1.87"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("D33").Value = "'21.11"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'157.26"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "'4.65"
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("D37").Value = "'6.19"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("D38").Value = "'26.76"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").Value = "'1.25"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "'0.0666"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "'3.99"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.53"
$ws.Range("E42").Value = "  +11.03%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.683"
$ws.Range("E43").Value = "  +4.85%  "
$ws.Range("D44").Value = "3.143.36"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "'36.87"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0262"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.281.57"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +7.22%  "
$ws.Range("D50").Value = "'20.67"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'6.09"
$ws.Range("E51").Value = "  +3.17%  "
